# This workbook holds weekly "Haba" (broad bean) price records for the
# "Macroferia Regional de Talca" market, one row per week, with the most
# recent week's data entered at the top of the data block (row 66) and
# older weeks pushed down below it.
#
# A new weekly record is being added at the top (row 66, date 2022-09-21 /
# serial 44825). Every existing record from row 66 down to row 90 shifts
# down by one row (66->67, 67->68, ... 89->90), and the record that used to
# be in row 90 now lands in a brand-new row 91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room at the bottom -------------------------------------
# Row 91 does not exist yet. Give its date cell (D91) the same number
# format as the rest of the date column *before* writing any value into it,
# so the engine doesn't invent a brand-new "general date" style for the
# previously-empty cell.
$ws.Range("D91").NumberFormat = $ws.Range("D90").NumberFormat()

# --- Step 2: shift every record in rows 66-90 down by one row ------------
# Walk from the bottom up so we never clobber a row before reading it.
for ($r = 90; $r -ge 66; $r--) {
    $srcRange = "A" + $r + ":R" + $r
    $dstRange = "A" + ($r + 1) + ":R" + ($r + 1)
    $rowVals = $ws.Range($srcRange).Value()
    $ws.Range($dstRange).Value = $rowVals
}

# --- Step 3: write the new record's date into row 66 ----------------------
# Everything else in row 66 (market, category, volume, prices, origin, etc.)
# already holds the correct values for the new record after the shift above
# (it's the same data that used to describe the prior top record); only the
# date actually changes for the newly-added week.
$ws.Range("D66").Value = 44825
